# Append new Lancers listings captured at 2026-01-12 12:42:28 JST.
# Row 2 keeps its original content (AI tech-lead listing) but gets the
# refreshed capture timestamp. A brand-new listing (Shopee API tool) is
# inserted as row 3, pushing the former row-3 (Zapier) listing down to
# row 4. Two more brand-new listings land in rows 5-6, the former row-4
# (HITOON) listing moves to row 7, and one final new listing becomes
# row 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2026-01-12 12:42:28"

# The engine's Hyperlinks.Delete() call clears every hyperlink on the
# whole sheet (not just the target cell's), so wipe them all up front
# and re-add every row's link afterwards, in final row order, below.
$ws.Range("F2").Hyperlinks.Delete()

# --- Row 2 : existing AI tech-lead listing, timestamp refreshed only ---
$ws.Range("A2").Value = $timestamp
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5423720")
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5423720"
$ws.Range("F2").Style = "Hyperlink"

# --- Row 3 (NEW) : Shopee API tool ---
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "Shopee APIを使用した「商品動画の一括紐付けツール」の開発依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5469483")
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5469483"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("G3").Value = 308
$ws.Range("H3").Value = "🔥API ◆ツール,開発"

# --- Row 4 : former row-3 listing (Zapier), timestamp refreshed ---
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "【Zapier設定のみ!作業時間~1時間】スプレッドシート・Gドライブ自動化構築(設計済)"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5469379")
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5469379"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("G4").Value = 255
$ws.Range("H4").Value = "🔥API ◆自動化"

# --- Row 5 (NEW) : Web app / management system ---
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "【フルリモート可】Webアプリ開発経験者募集!経営管理システムの開発"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5469430")
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5469430"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("G5").Value = 150
$ws.Range("H5").Value = "◆開発 ◇アプリ"

# --- Row 6 (NEW) : Java/Javascript engineer ---
$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "【未経験相談可能】JavaまたはJavascriptエンジニアを募集!"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5469522")
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5469522"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("G6").Value = 85
$ws.Range("H6").Value = "★Java"

# --- Row 7 : former row-4 listing (HITOON), timestamp refreshed ---
$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "フロント実装済み!音楽権利マーケットプレイス「HITOON」のバックエンド・決済機能実装"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5469298")
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5469298"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("G7").Value = 18
# (no H7 value - this listing has no skill-summary tag, same as before)

# --- Row 8 (NEW) : Microsoft Access tool ---
$ws.Range("A8").Value = $timestamp
$ws.Range("B8").Value = "【急募】Microsoft Accessで物流納品先別仕分けリスト作成"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5469531")
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5469531"
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("G8").Value = 10

# --- Column H a touch wider to fit the longer "◆ツール,開発" tag ---
$ws.Columns.Item(8).ColumnWidth = 13.14
